$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 29 - this shifts the existing rows 29-34 down to 30-35,
# copying row 29's formatting (incl. the date style on column D) to the new row.
$ws.Rows.Item(29).Insert()

# Populate the newly inserted row 29 with the new weekly record.
$ws.Range("A29").Value = 2
$ws.Range("B29").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C29").Value = "Coquimbo"
$ws.Range("D29").Value = 44468
$ws.Range("E29").Value = 4
$ws.Range("F29").Value = 100112022
$ws.Range("G29").Value = "Arveja Verde"
$ws.Range("H29").Value = "Perfection"
$ws.Range("I29").Value = "Primera"
$ws.Range("J29").Value = 500
$ws.Range("K29").Value = 23000
$ws.Range("L29").Value = 25000
$ws.Range("M29").Value = 24000
$ws.Range("N29").Value = '$/malla 25 kilos'
$ws.Range("O29").Value = "Provincia de Limarí"
$ws.Range("P29").Value = 960
$ws.Range("Q29").Value = 25
$ws.Range("R29").Value = "Hortaliza"
